$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H121").Value = 2033.6
$ws.Range("I121").Value = 0
$ws.Range("K121").Value = 0
$ws.Range("M121").Value = ""
$ws.Range("H132").Value = 1132.2307
$ws.Range("I132").Value = 1149.4783
$ws.Range("K132").Value = 3448.4349
$ws.Range("M132").Value = -918.4349000000002
$ws.Range("H137").Value = 1728.3182
$ws.Range("I137").Value = 1459
$ws.Range("J137").Value = 1997.6364
$ws.Range("K137").Value = 4377
$ws.Range("L137").Value = 5992.9092
$ws.Range("M137").Value = -1827
$ws.Range("N137").Value = -11092.9092
$ws.Range("H138").Value = 3306.5
$ws.Range("J138").Value = 2195.353
$ws.Range("L138").Value = 6586.059
$ws.Range("N138").Value = -16866.059

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1551692.4
$ws.Range("I2").Value = 2115317
$ws.Range("K2").Value = 2115317
$ws.Range("M2").Value = -2115204
$ws.Range("H23").Value = 45003
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").Value = ""
$ws.Range("H32").Value = 2690.1282
$ws.Range("J32").Value = 5329.8667
$ws.Range("L32").Value = 5329.8667
$ws.Range("N32").Value = -5903.8667
$ws.Range("H61").Value = 2801.5833
$ws.Range("I61").Value = 955.8570999999999
$ws.Range("J61").Value = 5385.6
$ws.Range("K61").Value = 955.8570999999999
$ws.Range("L61").Value = 5385.6
$ws.Range("M61").Value = -743.8570999999999
$ws.Range("N61").Value = -5809.6
$ws.Range("H63").Value = 9999.5
$ws.Range("I63").Value = 9999.5
$ws.Range("K63").Value = 9999.5
$ws.Range("M63").Value = -9313.5
$ws.Range("H66").Value = 9999.5
$ws.Range("I66").Value = 9999.5
$ws.Range("K66").Value = 49997.5
$ws.Range("M66").Value = -46565.5
$ws.Range("H74").Value = 1122.08
$ws.Range("I74").Value = 732.94116
$ws.Range("J74").Value = 1949
$ws.Range("K74").Value = 732.94116
$ws.Range("L74").Value = 1949
$ws.Range("M74").Value = 141.05884
$ws.Range("N74").Value = -3697
$ws.Range("H77").Value = 1122.08
$ws.Range("I77").Value = 732.94116
$ws.Range("J77").Value = 1949
$ws.Range("K77").Value = 3664.7058
$ws.Range("L77").Value = 9745
$ws.Range("M77").Value = 703.2942000000003
$ws.Range("N77").Value = -18481
$ws.Range("H116").Value = 1551692.4
$ws.Range("I116").Value = 2115317
$ws.Range("K116").Value = 2115317
$ws.Range("M116").Value = -2113023
$ws.Range("H136").Value = 2801.5833
$ws.Range("I136").Value = 955.8570999999999
$ws.Range("J136").Value = 5385.6
$ws.Range("K136").Value = 2867.5713
$ws.Range("L136").Value = 16156.8
$ws.Range("M136").Value = -317.5712999999996
$ws.Range("N136").Value = -21256.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1551692.4
$ws.Range("I3").Value = 2115317
$ws.Range("K3").Value = 2115317
$ws.Range("M3").Value = -2115203
$ws.Range("H80").Value = 9371.909
$ws.Range("I80").Value = 280
$ws.Range("J80").Value = 10281.1
$ws.Range("K80").Value = 280
$ws.Range("L80").Value = 10281.1
$ws.Range("M80").Value = 718
$ws.Range("N80").Value = -12277.1
$ws.Range("H83").Value = 9371.909
$ws.Range("I83").Value = 280
$ws.Range("J83").Value = 10281.1
$ws.Range("K83").Value = 1400
$ws.Range("L83").Value = 51405.5
$ws.Range("M83").Value = 3592
$ws.Range("N83").Value = -61389.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1471.9608
$ws.Range("I31").Value = 733.3889
$ws.Range("J31").Value = 1874.8182
$ws.Range("K31").Value = 733.3889
$ws.Range("L31").Value = 1874.8182
$ws.Range("M31").Value = -438.3889
$ws.Range("N31").Value = -2464.8182
$ws.Range("H34").Value = 1471.9608
$ws.Range("I34").Value = 733.3889
$ws.Range("J34").Value = 1874.8182
$ws.Range("K34").Value = 733.3889
$ws.Range("L34").Value = 1874.8182
$ws.Range("M34").Value = -531.3889
$ws.Range("N34").Value = -2278.8182
$ws.Range("H58").Value = 2416569.5
$ws.Range("I58").Value = 3953502.2
$ws.Range("K58").Value = 3953502.2
$ws.Range("M58").Value = -3953299.2
$ws.Range("H86").Value = 83335630
$ws.Range("I86").Value = 111113430
$ws.Range("K86").Value = 111113430
$ws.Range("M86").Value = -111112307
$ws.Range("H89").Value = 83335630
$ws.Range("I89").Value = 111113430
$ws.Range("K89").Value = 555567150
$ws.Range("M89").Value = -555561534
$ws.Range("H94").Value = 841.5714
$ws.Range("I94").Value = 738.5
$ws.Range("J94").Value = 1099.25
$ws.Range("K94").Value = 738.5
$ws.Range("L94").Value = 1099.25
$ws.Range("M94").Value = -287.5
$ws.Range("N94").Value = -2001.25
$ws.Range("H99").Value = 2981.111
$ws.Range("I99").Value = 2832.8572
$ws.Range("K99").Value = 2832.8572
$ws.Range("M99").Value = -1334.8572
$ws.Range("H126").Value = 2981.111
$ws.Range("I126").Value = 2832.8572
$ws.Range("K126").Value = 8498.571599999999
$ws.Range("M126").Value = -6028.571599999999
$ws.Range("H136").Value = 2416569.5
$ws.Range("I136").Value = 3953502.2
$ws.Range("K136").Value = 11860506.6
$ws.Range("M136").Value = -11857956.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2022733.6
$ws.Range("I4").Value = 3417883.8
$ws.Range("J4").Value = 976371
$ws.Range("K4").Value = 10253651.4
$ws.Range("L4").Value = 2929113
$ws.Range("M4").Value = -10253539.4
$ws.Range("N4").Value = -2929337
$ws.Range("H46").Value = 2249.5
$ws.Range("H68").Value = 2592.7407
$ws.Range("J68").Value = 2680.48
$ws.Range("L68").Value = 8041.440000000001
$ws.Range("N68").Value = -9663.440000000001
$ws.Range("H71").Value = 2592.7407
$ws.Range("J71").Value = 2680.48
$ws.Range("L71").Value = 24124.32
$ws.Range("N71").Value = -32236.32
$ws.Range("H107").Value = 1748.2273
$ws.Range("J107").Value = 2057.6428
$ws.Range("L107").Value = 6172.928400000001
$ws.Range("N107").Value = -10012.9284

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 879999.5
$ws.Range("J3").Value = 6666
$ws.Range("L3").Value = 6666
$ws.Range("N3").Value = -6898
$ws.Range("H11").Value = 6324869
$ws.Range("I11").Value = 2582236.5
$ws.Range("J11").Value = 11794870
$ws.Range("K11").Value = 2582236.5
$ws.Range("L11").Value = 11794870
$ws.Range("M11").Value = -2582097.5
$ws.Range("N11").Value = -11795148
$ws.Range("H122").Value = 2403.6667
$ws.Range("I122").Value = 1617.8334
$ws.Range("K122").Value = 4853.5002
$ws.Range("M122").Value = -2403.5002
$ws.Range("H126").Value = 2097465
$ws.Range("I126").Value = 2648406.8
$ws.Range("K126").Value = 7945220.399999999
$ws.Range("M126").Value = -7942750.399999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 15000
$ws.Range("J20").Value = 15000
$ws.Range("L20").Value = 15000
$ws.Range("N20").Value = -15452
$ws.Range("H40").Value = 6709.7
$ws.Range("I40").Value = 3683
$ws.Range("J40").Value = 11249.75
$ws.Range("K40").Value = 3683
$ws.Range("L40").Value = 11249.75
$ws.Range("M40").Value = -3547
$ws.Range("N40").Value = -11521.75
$ws.Range("H61").Value = 2289.4211
$ws.Range("I61").Value = 1833.3334
$ws.Range("J61").Value = 3999.75
$ws.Range("K61").Value = 1833.3334
$ws.Range("L61").Value = 3999.75
$ws.Range("M61").Value = -1631.3334
$ws.Range("N61").Value = -4403.75
$ws.Range("H63").Value = 16997
$ws.Range("J63").Value = 16997
$ws.Range("L63").Value = 16997
$ws.Range("N63").Value = -18495
$ws.Range("H66").Value = 16997
$ws.Range("J66").Value = 16997
$ws.Range("L66").Value = 50991
$ws.Range("N66").Value = -58479
$ws.Range("H113").Value = 2289.4211
$ws.Range("I113").Value = 1833.3334
$ws.Range("J113").Value = 3999.75
$ws.Range("K113").Value = 1833.3334
$ws.Range("L113").Value = 3999.75
$ws.Range("M113").Value = 336.6666
$ws.Range("N113").Value = -8339.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 13229567
$ws.Range("I136").Value = 34724860
$ws.Range("J136").Value = 1694.0385
$ws.Range("K136").Value = 104174580
$ws.Range("L136").Value = 5082.1155
$ws.Range("M136").Value = -104172030
$ws.Range("N136").Value = -10182.1155
